$wb = $excel.ActiveWorkbook

$activeBefore = $wb.ActiveSheet.Name

$ws2 = $wb.Worksheets.Item("CDCF-PMpPDOU")
$ws2.Range("B2").Formula = "=1/1.60934*10^12"
$ws2.Range("B10").Select()

$ws3 = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$ws3.Range("B2").Formula = "=1/1.60934*10^12"
$ws3.Range("B25").Select()

# Restore the originally active sheet so the saved workbook's active tab
# (and tabSelected flag) is unchanged by our selection calls above.
$wb.Worksheets.Item($activeBefore).Activate()
